# Update the "Last Updated" timestamp on the Metadata sheet.
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 10:39 AM"

# The "Stock List" sheet gained a brand-new top row (CAPTRU-RE1) and lost
# its previous last row (TRAVELFOOD) - i.e. every existing data row shifted
# down by one and the old final row fell off the bottom.
$ws = $wb.Worksheets.Item("Stock List")

# Insert a fresh row right under the header; this pushes all existing
# data rows (old row 2 .. old row 76) down by one automatically.
$ws.Rows.Item(2).Insert()

# Row insertion in Excel inherits formatting from the row above (the bold,
# centered header row here); the original data rows carry no explicit
# style, so strip whatever got copied in before writing values.
$ws.Range("A2:H2").ClearFormats()

# Populate the newly inserted row 2 with the new stock's data.
$ws.Range("A2").Value = [char]0x1F4CB
$ws.Range("B2").Value = "CAPTRU-RE1"
$ws.Range("C2").Value = "CAPTRU-RE1"
$ws.Range("D2").Value = 5.67
$ws.Range("E2").Value = -11.9565
$ws.Range("F2").Value = "N/A"
$ws.Range("G2").Value = "N/A"
$ws.Range("H2").Value = 0

# The insert above shifted the old last data row (TRAVELFOOD) down to
# row 77; remove it so the sheet keeps its original 76-row extent.
$ws.Rows.Item(77).Delete()
